$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (columns D, H, I, J, K, L, M, P)
$data = @{
    2  = @{ D = 44699; H = "Sin especificar"; I = "Primera"; J = 50;  K = 9000; L = 9500; M = 9250; P = 9250 }
    3  = @{ D = 44309; H = "Sin especificar"; I = "Primera"; J = 50;  K = 8000; L = 9000; M = 8500; P = 8500 }
    4  = @{ D = 44804; H = "Sin especificar"; I = "Primera"; J = 60;  K = 5500; L = 6000; M = 5750; P = 5750 }
    5  = @{ D = 44410; H = "Sin especificar"; I = "Primera"; J = 100; K = 5500; L = 6000; M = 5750; P = 5750 }
    6  = @{ D = 44259; H = "Sin especificar"; I = "Primera"; J = 80;  K = 4000; L = 4500; M = 4250; P = 4250 }
    7  = @{ D = 44539; H = "Americana (o)";   I = "Primera"; J = 160; K = 6500; L = 7000; M = 6750; P = 6750 }
    8  = @{ D = 44497; H = "Sin especificar"; I = "Primera"; J = 160; K = 5000; L = 6000; M = 5500; P = 5500 }
    9  = @{ D = 44281; H = "Sin especificar"; I = "Primera"; J = 100; K = 5000; L = 6000; M = 5500; P = 5500 }
    10 = @{ D = 44636; H = "Americana (o)";   I = "Primera"; J = 60;  K = 8000; L = 9000; M = 8500; P = 8500 }
    11 = @{ D = 44764; H = "Americana (o)";   I = "Primera"; J = 100; K = 7000; L = 8000; M = 7500; P = 7500 }
    12 = @{ D = 44263; H = "Sin especificar"; I = "Primera"; J = 100; K = 7000; L = 8000; M = 7500; P = 7500 }
    13 = @{ D = 44789; H = "Sin especificar"; I = "Primera"; J = 80;  K = 5000; L = 6000; M = 5500; P = 5500 }
    14 = @{ D = 44253; H = "Americana (o)";   I = "Segunda"; J = 100; K = 4000; L = 4500; M = 4250; P = 4250 }
    15 = @{ D = 44371; H = "Sin especificar"; I = "Primera"; J = 80;  K = 7000; L = 8000; M = 7375; P = 7375 }
    16 = @{ D = 44575; H = "Sin especificar"; I = "Primera"; J = 160; K = 6500; L = 7000; M = 6750; P = 6750 }
    17 = @{ D = 44414; H = "Sin especificar"; I = "Primera"; J = 100; K = 6000; L = 7000; M = 6500; P = 6500 }
    18 = @{ D = 44559; H = "Americana (o)";   I = "Primera"; J = 100; K = 5000; L = 6000; M = 5500; P = 5500 }
    19 = @{ D = 44945; H = "Sin especificar"; I = "Primera"; J = 45;  K = 6000; L = 7000; M = 6444; P = 6444 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals.D   # D - Fecha
    $ws.Cells.Item($row, 8).Value = $vals.H   # H - Variedad
    $ws.Cells.Item($row, 9).Value = $vals.I   # I - Calidad
    $ws.Cells.Item($row, 10).Value = $vals.J  # J - Volumen
    $ws.Cells.Item($row, 11).Value = $vals.K  # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals.L  # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals.M  # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals.P  # P - Precio $/Kg
}

$wb.Save()
